# Generate Report for handoff
#
# The previous handoff attempt (uuid 9c1df6c6-e178-419a-820c-5627076f4fcf) is
# superseded by a new report run (uuid 77cb3479-5038-4779-9011-473fd00eda82).
# Because the handoff transform for this file failed, the per-language rows
# get their "last handoff" bookkeeping cleared out and the row is marked to
# be ignored rather than included in the next handoff.

$wb = $excel.ActiveWorkbook

$newFile = "77cb3479-5038-4779-9011-473fd00eda82.md"
$newStatus = "Handoff transform failed"
$zeroDate = "0001-01-01 00:00:00"

# --- Overview sheet: just the file name + status roll up ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = $newFile
$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus

# --- zh-cn sheet ---
$wsZh = $wb.Worksheets.Item("zh-cn")
# Drop the stale "Latest Handoff File" hyperlink (C2) along with the rest of
# the sheet's hyperlinks; the ones that still apply (A2/A3) point at cells
# whose text is rewritten below.
$wsZh.Hyperlinks.Delete()
$wsZh.Range("A2").Value = $newFile
$wsZh.Range("B2").Value = $newStatus
$wsZh.Range("C2").Clear()
$wsZh.Range("D2").Value = $zeroDate
$wsZh.Range("H2").Value = "Ignored"

# --- de-de sheet ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Hyperlinks.Delete()
$wsDe.Range("A2").Value = $newFile
$wsDe.Range("B2").Value = $newStatus
$wsDe.Range("C2").Clear()
$wsDe.Range("D2").Value = $zeroDate
$wsDe.Range("H2").Value = "Ignored"
